$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1646.3462
$ws.Range("I28").Value = 474.0625
$ws.Range("J28").Value = 3522
$ws.Range("K28").Value = 474.0625
$ws.Range("L28").Value = 3522
$ws.Range("M28").Value = 10.9375
$ws.Range("N28").Value = -4492
$ws.Range("H33").Value = 608.9167
$ws.Range("I33").Value = 243.125
$ws.Range("K33").Value = 243.125
$ws.Range("M33").Value = -14.125
$ws.Range("H70").Value = 2602.7856
$ws.Range("I70").Value = 2380
$ws.Range("J70").Value = 3003.8
$ws.Range("K70").Value = 7140
$ws.Range("L70").Value = 9011.400000000001
$ws.Range("M70").Value = -6870
$ws.Range("N70").Value = -9551.400000000001
$ws.Range("H73").Value = 2602.7856
$ws.Range("I73").Value = 2380
$ws.Range("J73").Value = 3003.8
$ws.Range("K73").Value = 7140
$ws.Range("L73").Value = 9011.400000000001
$ws.Range("M73").Value = -6204
$ws.Range("N73").Value = -10883.4
$ws.Range("H86").Value = 6837.5
$ws.Range("I86").Value = 19000
$ws.Range("J86").Value = 5731.8184
$ws.Range("K86").Value = 19000
$ws.Range("L86").Value = 5731.8184
$ws.Range("M86").Value = -17877
$ws.Range("N86").Value = -7977.8184
$ws.Range("H89").Value = 6837.5
$ws.Range("I89").Value = 19000
$ws.Range("J89").Value = 5731.8184
$ws.Range("K89").Value = 95000
$ws.Range("L89").Value = 28659.092
$ws.Range("M89").Value = -89384
$ws.Range("N89").Value = -39891.092
$ws.Range("H111").Value = 5558713
$ws.Range("I111").Value = 9261351
$ws.Range("J111").Value = 4756.125
$ws.Range("K111").Value = 27784053
$ws.Range("L111").Value = 14268.375
$ws.Range("M111").Value = -27780986
$ws.Range("N111").Value = -20402.375
$ws.Range("H132").Value = 2004.7931
$ws.Range("I132").Value = 1991.3091
$ws.Range("J132").Value = 2252
$ws.Range("K132").Value = 5973.927299999999
$ws.Range("L132").Value = 6756
$ws.Range("M132").Value = -3443.927299999999
$ws.Range("N132").Value = -11816
$ws.Range("H137").Value = 43217.387
$ws.Range("J137").Value = 5267.2
$ws.Range("L137").Value = 15801.6
$ws.Range("N137").Value = -20901.6
$ws.Range("H138").Value = 3564.5144
$ws.Range("I138").Value = 2856.8462
$ws.Range("J138").Value = 3725.9124
$ws.Range("K138").Value = 8570.5386
$ws.Range("L138").Value = 11177.7372
$ws.Range("M138").Value = -3430.5386
$ws.Range("N138").Value = -21457.7372

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9390.531999999999
$ws.Range("I32").Value = 7790.4365
$ws.Range("K32").Value = 7790.4365
$ws.Range("M32").Value = -7503.4365
$ws.Range("H61").Value = 5457.636
$ws.Range("I61").Value = 5056.2104
$ws.Range("K61").Value = 5056.2104
$ws.Range("M61").Value = -4844.2104
$ws.Range("H74").Value = 36315.297
$ws.Range("I74").Value = 2749.5293
$ws.Range("J74").Value = 93377.10000000001
$ws.Range("K74").Value = 2749.5293
$ws.Range("L74").Value = 93377.10000000001
$ws.Range("M74").Value = -1875.5293
$ws.Range("N74").Value = -95125.10000000001
$ws.Range("H77").Value = 36315.297
$ws.Range("I77").Value = 2749.5293
$ws.Range("J77").Value = 93377.10000000001
$ws.Range("K77").Value = 13747.6465
$ws.Range("L77").Value = 466885.5
$ws.Range("M77").Value = -9379.646500000001
$ws.Range("N77").Value = -475621.5
$ws.Range("H88").Value = 3019.125
$ws.Range("J88").Value = 1491.5
$ws.Range("L88").Value = 1491.5
$ws.Range("N88").Value = -2303.5
$ws.Range("H91").Value = 3019.125
$ws.Range("J91").Value = 1491.5
$ws.Range("L91").Value = 1491.5
$ws.Range("N91").Value = -4299.5
$ws.Range("H132").Value = 36366.92
$ws.Range("I132").Value = 1666.2727
$ws.Range("K132").Value = 4998.8181
$ws.Range("M132").Value = -2468.8181
$ws.Range("H136").Value = 5457.636
$ws.Range("I136").Value = 5056.2104
$ws.Range("K136").Value = 15168.6312
$ws.Range("M136").Value = -12618.6312

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 35000
$ws.Range("J38").Value = 35000
$ws.Range("L38").Value = 35000
$ws.Range("N38").Value = -35832
$ws.Range("H80").Value = 740.2308
$ws.Range("I80").Value = 830.8
$ws.Range("J80").Value = 683.625
$ws.Range("K80").Value = 830.8
$ws.Range("L80").Value = 683.625
$ws.Range("M80").Value = 167.2
$ws.Range("N80").Value = -2679.625
$ws.Range("H83").Value = 740.2308
$ws.Range("I83").Value = 830.8
$ws.Range("J83").Value = 683.625
$ws.Range("K83").Value = 4154
$ws.Range("L83").Value = 3418.125
$ws.Range("M83").Value = 838
$ws.Range("N83").Value = -13402.125
$ws.Range("H134").Value = 6167.7085
$ws.Range("I134").Value = 1955.909
$ws.Range("K134").Value = 5867.727000000001
$ws.Range("M134").Value = -3332.727000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20795.143
$ws.Range("I31").Value = 3179.5334
$ws.Range("J31").Value = 41120.848
$ws.Range("K31").Value = 3179.5334
$ws.Range("L31").Value = 41120.848
$ws.Range("M31").Value = -2884.5334
$ws.Range("N31").Value = -41710.848
$ws.Range("H34").Value = 20795.143
$ws.Range("I34").Value = 3179.5334
$ws.Range("J34").Value = 41120.848
$ws.Range("K34").Value = 3179.5334
$ws.Range("L34").Value = 41120.848
$ws.Range("M34").Value = -2977.5334
$ws.Range("N34").Value = -41524.848
$ws.Range("H50").Value = 4088.889
$ws.Range("J50").Value = 4088.889
$ws.Range("L50").Value = 4088.889
$ws.Range("N50").Value = -5338.889
$ws.Range("H99").Value = 4717.636
$ws.Range("I99").Value = 4237.375
$ws.Range("J99").Value = 5998.3335
$ws.Range("K99").Value = 4237.375
$ws.Range("L99").Value = 5998.3335
$ws.Range("M99").Value = -2739.375
$ws.Range("N99").Value = -8994.333500000001
$ws.Range("H121").Value = 44498.75
$ws.Range("J121").Value = 44498.75
$ws.Range("L121").Value = 44498.75
$ws.Range("N121").Value = -47118.75
$ws.Range("H126").Value = 4717.636
$ws.Range("I126").Value = 4237.375
$ws.Range("J126").Value = 5998.3335
$ws.Range("K126").Value = 12712.125
$ws.Range("L126").Value = 17995.0005
$ws.Range("M126").Value = -10242.125
$ws.Range("N126").Value = -22935.0005
$ws.Range("H132").Value = 38168.188
$ws.Range("I132").Value = 26624.098
$ws.Range("J132").Value = 77610.5
$ws.Range("K132").Value = 79872.29400000001
$ws.Range("L132").Value = 232831.5
$ws.Range("M132").Value = -77342.29400000001
$ws.Range("N132").Value = -237891.5
$ws.Range("H141").Value = 143463.92
$ws.Range("J141").Value = 143463.92
$ws.Range("L141").Value = 143463.92
$ws.Range("N141").Value = -153823.92

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 4667.6665
$ws.Range("I7").Value = 4500.5
$ws.Range("K7").Value = 13501.5
$ws.Range("M7").Value = -13389.5
$ws.Range("H68").Value = 1654.2
$ws.Range("J68").Value = 2999.5
$ws.Range("L68").Value = 8998.5
$ws.Range("N68").Value = -10620.5
$ws.Range("H71").Value = 1654.2
$ws.Range("J71").Value = 2999.5
$ws.Range("L71").Value = 26995.5
$ws.Range("N71").Value = -35107.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1142620.6
$ws.Range("I80").Value = 1628408.4
$ws.Range("J80").Value = 215207.73
$ws.Range("K80").Value = 1628408.4
$ws.Range("L80").Value = 215207.73
$ws.Range("M80").Value = -1627410.4
$ws.Range("N80").Value = -217203.73
$ws.Range("H83").Value = 1142620.6
$ws.Range("I83").Value = 1628408.4
$ws.Range("J83").Value = 215207.73
$ws.Range("K83").Value = 8142042
$ws.Range("L83").Value = 1076038.65
$ws.Range("M83").Value = -8137050
$ws.Range("N83").Value = -1086022.65
$ws.Range("H102").Value = 5191456
$ws.Range("I102").Value = 7408757
$ws.Range("J102").Value = 2167863.8
$ws.Range("K102").Value = 7408757
$ws.Range("L102").Value = 2167863.8
$ws.Range("M102").Value = -7407135
$ws.Range("N102").Value = -2171107.8
$ws.Range("H107").Value = 11717.444
$ws.Range("I107").Value = 17232.834
$ws.Range("J107").Value = 686.6667
$ws.Range("K107").Value = 17232.834
$ws.Range("L107").Value = 686.6667
$ws.Range("M107").Value = -15312.834
$ws.Range("N107").Value = -4526.6667
$ws.Range("H132").Value = 2872.9207
$ws.Range("I132").Value = 2594.44
$ws.Range("J132").Value = 3944
$ws.Range("K132").Value = 7783.32
$ws.Range("L132").Value = 11832
$ws.Range("M132").Value = -5253.32
$ws.Range("N132").Value = -16892

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 416.66666
$ws.Range("I16").Value = 416.66666
$ws.Range("K16").Value = 416.66666
$ws.Range("M16").Value = -246.66666
$ws.Range("H61").Value = 7409000
$ws.Range("I61").Value = 7938071.5
$ws.Range("K61").Value = 7938071.5
$ws.Range("M61").Value = -7937869.5
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H113").Value = 7409000
$ws.Range("I113").Value = 7938071.5
$ws.Range("K113").Value = 7938071.5
$ws.Range("M113").Value = -7935901.5
$ws.Range("H122").Value = 5064.7295
$ws.Range("J122").Value = 6981.4546
$ws.Range("L122").Value = 20944.3638
$ws.Range("N122").Value = -25844.3638
$ws.Range("H132").Value = 11821.3
$ws.Range("I132").Value = 12878.883
$ws.Range("J132").Value = 5828.3335
$ws.Range("K132").Value = 38636.649
$ws.Range("L132").Value = 17485.0005
$ws.Range("M132").Value = -36106.649
$ws.Range("N132").Value = -22545.0005

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 6415544
$ws.Range("J81").Value = 9527
$ws.Range("L81").Value = 19054
$ws.Range("N81").Value = -21176
$ws.Range("H84").Value = 6415544
$ws.Range("J84").Value = 9527
$ws.Range("L84").Value = 95270
$ws.Range("N84").Value = -105878
